# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" worksheet (fund-holdings detail) right after the
# "总计" (totals) sheet, and updates the "总计" sheet with the corresponding
# summary row, shifting the previously-newest rows down one position.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) summary sheet: insert the 2022-Q4 row at the
#    top of the data and push the existing rows down by one.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Give the new last data row (A5) the same numeric-index style as A4 before
# writing into it.
$totals.Range("A4").Copy()
$totals.Range("A5").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 22
$totals.Range("D2").Value = 1.13

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 7
$totals.Range("D3").Value = 0.44

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2021-Q4"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 2.85

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q3"
$totals.Range("C5").Value = 4
$totals.Range("D5").Value = 0

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totals)
$q4.Name = "2022-Q4"

# Pull the header / index-column formatting (bold, centred, bordered style)
# from the totals sheet so the new sheet matches its siblings.
$totals.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$totals.Range("A2").Copy()
$q4.Range("A2:A23").PasteSpecial(-4122)

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B (基金代码), C (基金名称), D (基金规模), E (股票总仓位),
# F (仓位占比) and (for almost every row) G (持有市值) are stored as plain
# text in the source data (e.g. fund codes keep their leading zeros,
# percentages/market values keep trailing zeros) — force text formatting
# before writing so the values round-trip exactly. The couple of G cells
# that hold a genuine 0 are reset to General further below.
$q4.Range("B2:G23").NumberFormat = "@"

$rows = @"
0	001113	南方大数据100指数A	15.48	89.97	2.19	0.3390	S	4
1	501201	红土创新科技创新 3 年封闭	3.82	92.71	4.14	0.1581	S	8
2	006265	红土创新新科技股票	2.90	93.37	4.84	0.1404	S	8
3	009467	红土创新科技创新3个月定开混合A	1.62	88.86	5.34	0.0865	S	4
4	002707	摩根士丹利华鑫科技领先灵活配置混合A	1.76	92.23	4.84	0.0852	S	3
5	015005	中邮能源革新混合C	1.29	93.71	4.89	0.0631	S	6
6	880007	招商资管智远成长灵活配置混合A	1.49	90.93	3.49	0.0520	S	8
7	168401	红土创新转型精选灵活配置混合（LOF）	0.82	93.08	4.28	0.0351	S	7
8	012102	国寿安保低碳经济混合A	1.04	84.42	2.98	0.0310	S	9
9	013173	红土创新科技创新3个月定开混合C	0.43	88.86	5.34	0.0230	S	4
10	881007	招商资管智远成长灵活配置混合C	0.60	90.93	3.49	0.0209	S	8
11	001744	诺安进取回报灵活配置混合	0.59	69.55	3.53	0.0208	S	3
12	011729	工银聚享混合A	1.15	29.51	1.78	0.0205	S	3
13	165317	建信丰裕多策略灵活配置混合（LOF）	0.33	90.65	4.27	0.0141	S	7
14	011027	国寿安保稳弘混合A	0.62	24.35	2.08	0.0129	S	4
15	011028	国寿安保稳弘混合C	0.44	24.35	2.08	0.0092	S	4
16	015004	中邮能源革新混合A	0.12	93.71	4.89	0.0059	S	6
17	004344	南方大数据100指数C	0.21	89.97	2.19	0.0046	S	4
18	012103	国寿安保低碳经济混合C	0.13	84.42	2.98	0.0039	S	9
19	014871	摩根士丹利华鑫科技领先灵活配置混合C	0.05	92.23	4.84	0.0024	S	3
20	015407	国寿安保稳弘混合E	0.00	24.35	2.08	0	N	4
21	011730	工银聚享混合C	0.00	29.51	1.78	0	N	3
"@

$lines = $rows -split "`n"
$r = 2
foreach ($line in $lines) {
    $f = $line -split "`t"

    $q4.Cells.Item($r, 1).Value = [int]$f[0]
    $q4.Cells.Item($r, 2).Value = $f[1]
    $q4.Cells.Item($r, 3).Value = $f[2]
    $q4.Cells.Item($r, 4).Value = $f[3]
    $q4.Cells.Item($r, 5).Value = $f[4]
    $q4.Cells.Item($r, 6).Value = $f[5]

    if ($f[7] -eq "N") {
        $q4.Cells.Item($r, 7).NumberFormat = "General"
        $q4.Cells.Item($r, 7).Value = [double]$f[6]
    } else {
        $q4.Cells.Item($r, 7).Value = $f[6]
    }

    $q4.Cells.Item($r, 8).Value = [int]$f[8]

    $r = $r + 1
}

$q4.Range("A1").Select()
